# US-11840 [FIX] Return from Unit export: Changed template
#
# - "Sending date:" / "Delivery date:" labels renamed and the receipt-date
#   cell gets a bottom border + is merged across G11:H11.
# - Signature labels simplified to "Sent by:" / "Received by:" and the
#   row that holds them shrinks back to the default row height.
# - Selection / scroll position of the sheet view updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: date labels -----------------------------------------------
$ws.Range("C11").Value = "Creation date:"
$ws.Range("G11").Value = "Expected Receipt date:"

# G11:H11 becomes a merged cell with a thin bottom border (like an
# underline for the date to be filled in), matching C11:D11's "blank line".
$ws.Range("G11:H11").Merge()
$ws.Range("G11:H11").Borders.Item(9).LineStyle = 1
$ws.Range("G11:H11").Borders.Item(9).Weight = 2
$ws.Range("G11:H11").VerticalAlignment = -4108

# --- Row 24: signature labels -------------------------------------------
$ws.Range("B24").Value = "Sent by:"
$ws.Range("G24").Value = "Received by:"
$ws.Rows.Item(24).RowHeight = 15

# --- Sheet view: scroll position / selection -----------------------------
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 4 } catch {}
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("A4") } catch {}
$ws.Range("G25").Select()
